# Update the SAMHSA grant history table text to reflect FY 2012-2016
# (previously FY 2011-2016) in the descriptive cells above the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "This table shows the grant awards and award dollars SAMHSA made for FY 2012-2016. It is provided as a text alternative to the interactive chart on the SAMHSA page of this website."
$ws.Range("A7").Value = "Grant awards and award dollars SAMHSA made for FY 2012-2016."
